$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.833.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.943.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.17%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.563"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.03%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +1.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0890"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.97%  "

$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.409.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.942.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.987"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.942.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.70%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0987"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("E26").Value = "  +8.96%  "

$ws.Range("E27").Value = "  +3.13%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +17.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0452"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("E40").Value = "  +1.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.17%  "

$ws.Range("E42").Value = "  +1.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.62%  "

$ws.Range("E44").Value = "  -1.25%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.54%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.168.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.246"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("E50").Value = "  +10.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.935"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.96%  "
